# Development_Roadmap.xlsx edit script
# Summary of change (per commit message "Improved frame normalization, GUI adapted"):
#  - Insert a new task row "Frame brightness normalization" (done) right after the
#    "delete key" row, pushing all following rows down by one.
#  - Mark the three "Debayering" tasks (VideoReader / ImageReader / dark-flat import)
#    as "done" and update / simplify some of their discussion text.
#  - Re-prioritize "Documentation of I/O file formats" and "Documentation of
#    Debayering" from "Must have" to "Nice to have".
#  - Move the two embedded screenshots down by one row to keep them aligned with
#    their associated task rows.
#  - Update the window scroll position / selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 5 (shifts rows 5..27 down to 6..28, and the
#    used range grows from A1:G27 to A1:G28 automatically).
# ---------------------------------------------------------------------------
$ws.Rows(5).Insert()

# ---------------------------------------------------------------------------
# 2. Fill in the new row 5 with the "Frame brightness normalization" task.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value2 = 'Frame brightness normalization'
$ws.Range("B5").Value2 = "So far, only the related images (grayscale, Gauss, etc.) are normalized, such that the evaluation is not misled by brightness variations.`r`nNormalization should also be implemented for frame stacking. This will reduce blending artifacts."
$ws.Range("C5").Value2 = 'The derived frame types should not be changed at all. Only the brightnesses are stored as a single value per frame. This will reduce compute times and limit clipping defects. Also, normalization should be made optional by adding GUI parameters.'
$ws.Range("D5").Value2 = 'Rolf'
$ws.Range("E5").Value2 = 'Must have'
$ws.Range("F5").Value2 = '0.7.0'
$ws.Range("G5").Value2 = 'done'
$ws.Rows(5).RowHeight = 105

# ---------------------------------------------------------------------------
# 3. Row 7 ("Implementation of Debayering in VideoReader"): extend description,
#    mark as done, row grows taller.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value2 = 'The selection of debayering patterns is alredy implemented in the GUI. It is passed to the VideoReader initialization via the instance variable "self.debayer_pattern". All supported patterns should be recognized and implemented in the read_frame method. Automatic detection of bayer matrices is implemented as well as the automatic detection of RGB channel ordering.'
$ws.Range("G7").Value2 = 'done'
$ws.Rows(7).RowHeight = 90

# ---------------------------------------------------------------------------
# 4. Row 8 ("Implementation of Debayering in ImageReader"): replace the
#    discussion comment, mark as done.
# ---------------------------------------------------------------------------
$ws.Range("C8").Value2 = 'We assume that still image files are debayered already.'
$ws.Range("G8").Value2 = 'done'

# ---------------------------------------------------------------------------
# 5. Row 9 ("Debayering during dark / flat file import"): replace the long
#    clarification comment with the short resolution note, mark as done, and
#    shrink the row (no longer needs the huge height).
# ---------------------------------------------------------------------------
$ws.Range("C9").Value2 = 'We assume that still image files are debayered already.'
$ws.Range("G9").Value2 = 'done'
$ws.Rows(9).RowHeight = 75

# ---------------------------------------------------------------------------
# 6. Re-prioritize the two documentation rows (now rows 13 and 14) from
#    "Must have" to "Nice to have".
# ---------------------------------------------------------------------------
$ws.Range("E13").Value2 = 'Nice to have'
$ws.Range("E14").Value2 = 'Nice to have'

# ---------------------------------------------------------------------------
# 7. Move the two embedded pictures down by one row, preserving their
#    horizontal position, size and their offset within the anchor row.
# ---------------------------------------------------------------------------
$shapes = $ws.Shapes
$pic1 = $shapes.Item(1)
$pic2 = $shapes.Item(2)
$pic1.Top = $ws.Rows(18).Top + (628650 / 12700.0)
$pic2.Top = $ws.Rows(19).Top + (1343024 / 12700.0)

# ---------------------------------------------------------------------------
# 8. Update the view: scroll so row 4 is at the top, and select cell F14.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F14").Select()
